$p = $ppt.ActivePresentation

# --- 1. RR Data Manager hostname update (appears on the "Outbound" and
#        "Inbound" transfer workflow slides, shape "Rounded Rectangle 17") ---
$oldHost = "rr-datamgr01.hpc.psu.edu"
$newHost = "rr-datamgr.rr.hpc.psu.edu"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldHost) {
                    $tr.Text = $newHost
                }
            }
        }
    }
}

# --- 2. Refresh the auto date field ("datetimeFigureOut") everywhere it is
#        stamped: the slide master, every slide layout, and the notes
#        master all carry a "Date Placeholder" shape showing the date the
#        deck was last saved. ---
$oldDate = "10/9/24"
$newDate = "1/14/2025"

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $s = $shapes.Item($k)
        if ($s.HasTextFrame) {
            if ($s.TextFrame.HasText) {
                $tr = $s.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes
